# Weekly update: insert 2 new price records for Choclo (Vega Central Mapocho de
# Santiago) at the top of the data block starting at row 623. This shifts the
# existing rows 623-661 down to 625-663 (dimension grows from R661 to R663).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the start of the data block (before old row 623).
$ws.Rows("623:624").Insert()

# --- New row 623 ---
$ws.Range("A623").Value = 9
$ws.Range("B623").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C623").Value = "Metropolitana"
$ws.Range("D623").Value = 44931
$ws.Range("E623").Value = 13
$ws.Range("F623").Value = 100112024
$ws.Range("G623").Value = "Choclo"
$ws.Range("H623").Value = "Choclero"
$ws.Range("I623").Value = "Primera"
$ws.Range("J623").Value = 8800
$ws.Range("K623").Value = 350
$ws.Range("L623").Value = 380
$ws.Range("M623").Value = 365
$ws.Range("N623").Value = "$/unidad"
$ws.Range("O623").Value = "Provincia de Melipilla"
$ws.Range("P623").Value = 365
$ws.Range("Q623").Value = 1
$ws.Range("R623").Value = "Hortaliza"

# --- New row 624 ---
$ws.Range("A624").Value = 9
$ws.Range("B624").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C624").Value = "Metropolitana"
$ws.Range("D624").Value = 44931
$ws.Range("E624").Value = 13
$ws.Range("F624").Value = 100112024
$ws.Range("G624").Value = "Choclo"
$ws.Range("H624").Value = "Choclero"
$ws.Range("I624").Value = "Primera"
$ws.Range("J624").Value = 9700
$ws.Range("K624").Value = 350
$ws.Range("L624").Value = 380
$ws.Range("M624").Value = 365
$ws.Range("N624").Value = "$/unidad"
$ws.Range("O624").Value = "Región de O'Higgins"
$ws.Range("P624").Value = 365
$ws.Range("Q624").Value = 1
$ws.Range("R624").Value = "Hortaliza"
